# Apply the "gridOperator" config changes:
#  - config_gridNodes (sheet2): add an "operator" column (H) with "o1" for
#    the four electric grid nodes (rows 2-5); the heat node (row 6) is left
#    without an operator.
#  - config_actors (sheet4): add a new "gridoperator" actor row (id "o1").
#  - update the active sheet / selections to match the final state.

$wb = $excel.ActiveWorkbook

# ---- config_gridNodes: new "operator" column -----------------------------
$wsNodes = $wb.Worksheets.Item("config_gridNodes")

$wsNodes.Range("H1").Value = "operator"
$wsNodes.Range("H2").Value = "o1"
$wsNodes.Range("H3").Value = "o1"
$wsNodes.Range("H4").Value = "o1"
$wsNodes.Range("H5").Value = "o1"

# ---- config_actors: new gridoperator actor row ----------------------------
$wsActors = $wb.Worksheets.Item("config_actors")

$wsActors.Range("A25").Value = 23
$wsActors.Range("B25").Value = "gridoperator"
$wsActors.Range("C25").Value = "GRIDOPERATOR"
$wsActors.Range("D25").Value = "o1"

# ---- selections / active sheet --------------------------------------------
# Final active tab is config_actors; config_gridNodes keeps a leftover
# selection but is no longer the selected tab.
$wsNodes.Activate()
$wsNodes.Range("G14").Select()

$wsActors.Activate()
$wsActors.Range("B25").Select()
